# Update the "Project repo ..." line on the Contacts slide (slide 22) so
# that it reads:
#   "Project repo: https://github.com/0venoven/Peruse , contact us to be
#    added into the repo :D"
# with the URL rendered as an underlined hyperlink run (scheme color
# "hlink") that links to https://github.com/0venoven/Peruse - matching the
# style already used for the "Ivan"/"Yin Kit" github links above it.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(22)
$sh = $s.Shapes.Item(2)

# --- remember the shape's exact original position/size (in EMU) -----------
# Reading Shape.Left/Top/Width/Height back via COM truncates to a 32-bit
# float, which is not precise enough to reproduce the original EMU values
# after a round trip. Editing the run text below makes the auto-fit shape
# recompute its height, so we restore the original geometry afterwards
# using the EMU numbers straight from the source XML (not the lossy
# getters), nudged by a hair so the float truncation still rounds to the
# correct EMU integer.
$EMU_PER_PT = 12700.0
$EPS = 0.00005

$origLeftEmu   = 344500
$origTopEmu    = 3809850
$origWidthEmu  = 6015600
$origHeightEmu = 1154400

$tf = $sh.TextFrame
$tr = $tf.TextRange

# Paragraph 4 is the "Project repo is on GitHub ..." paragraph.
$para = $tr.Paragraphs(4, 1)
$start = $para.Start
$len = $para.Length

# Replace the whole paragraph with the new wording (select the exact
# existing span via Characters so the run isn't arbitrarily re-split).
$newText = "Project repo: https://github.com/0venoven/Peruse , contact us to be added into the repo :D"
$whole = $tr.Characters($start, $len)
$whole.Text = $newText

# Style+link the URL portion only.
$linkStart = $start + 14   # length of "Project repo: "
$linkLen = 34               # length of "https://github.com/0venoven/Peruse"
$linkRange = $tr.Characters($linkStart, $linkLen)
$linkRange.Font.Underline = $true
$linkRange.Font.Color.ObjectThemeColor = 11  # msoThemeColorHyperlink
$linkRange.ActionSettings.Item(1).Hyperlink.Address = "https://github.com/0venoven/Peruse"

# --- restore the shape's original geometry ---------------------------------
$sh.Left   = ($origLeftEmu   / $EMU_PER_PT) + $EPS
$sh.Top    = ($origTopEmu    / $EMU_PER_PT) + $EPS
$sh.Width  = ($origWidthEmu  / $EMU_PER_PT) + $EPS
$sh.Height = ($origHeightEmu / $EMU_PER_PT) + $EPS
